$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is a table of "field generator" rows (Field / Type / Generate /
# Pattern). Add a new "Sequence generator" row right after the existing
# "Auto number" row (row 20), reusing that row's look & feel, the same way
# a user would duplicate the last row of the table and then edit its text.
$ws.Range("A19:D19").Copy($ws.Range("A21:D21"))
$ws.Rows("21:21").RowHeight = 15.75

$ws.Range("A21").Value = "Sequence generator"
$ws.Range("B21").Value = "Sequence"
$ws.Range("C21").Value = "y"
$ws.Range("D21").Value = "ALP01_5_[SEQ(4,3,FAIL)] GEN"

# Leave the selection/scroll where a user would land after typing the new
# row: just past it, ready for the next entry.
$ws.Range("D22").Select()
$excel.ActiveWindow.ScrollRow = 7
